# Apply the Lich_Profits.xlsx value updates (ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets)
# via direct cell writes, mirroring the upstream scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 94
$ws.Range("H94").Value = 898.5
$ws.Range("I94").Value = 898.5
$ws.Range("K94").Value = 898.5
$ws.Range("M94").Value = -447.5

# Row 132
$ws.Range("H132").Value = 4205.3447
$ws.Range("I132").Value = 3664.1667
$ws.Range("K132").Value = 10992.5001
$ws.Range("M132").Value = -8462.500100000001

# Row 137
$ws.Range("H137").Value = 31209.564
$ws.Range("I137").Value = 43417
$ws.Range("K137").Value = 130251
$ws.Range("M137").Value = -127701

# Row 138
$ws.Range("H138").Value = 2244.318
$ws.Range("I138").Value = 819.5111000000001
$ws.Range("J138").Value = 3735.3953
$ws.Range("K138").Value = 2458.5333
$ws.Range("L138").Value = 11206.1859
$ws.Range("M138").Value = 2681.4667
$ws.Range("N138").Value = -21486.1859

# Row 141
$ws.Range("H141").Value = 2391.4375
$ws.Range("I141").Value = 2451.2222
$ws.Range("J141").Value = 2068.6
$ws.Range("K141").Value = 7353.6666
$ws.Range("L141").Value = 6205.799999999999
$ws.Range("M141").Value = -2173.6666
$ws.Range("N141").Value = -16565.8

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1853.9474
$ws.Range("I2").Value = 1706.7333
$ws.Range("K2").Value = 1706.7333
$ws.Range("M2").Value = -1593.7333

# Row 5
$ws.Range("H5").Value = 183.75
$ws.Range("J5").Value = 315
$ws.Range("L5").Value = 315
$ws.Range("N5").Value = -539

# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# Row 32
$ws.Range("H32").Value = 9387.745000000001
$ws.Range("I32").Value = 6427.8647
$ws.Range("K32").Value = 6427.8647
$ws.Range("M32").Value = -6140.8647

# Row 40
$ws.Range("H40").Value = 61249.25
$ws.Range("I40").Value = 48332.668
$ws.Range("J40").Value = 99999
$ws.Range("K40").Value = 48332.668
$ws.Range("L40").Value = 99999
$ws.Range("M40").Value = -48156.668
$ws.Range("N40").Value = -100351

# Row 74
$ws.Range("H74").Value = 105452.84
$ws.Range("I74").Value = 110755
$ws.Range("K74").Value = 110755
$ws.Range("M74").Value = -109881

# Row 77
$ws.Range("H77").Value = 105452.84
$ws.Range("I77").Value = 110755
$ws.Range("K77").Value = 553775
$ws.Range("M77").Value = -549407

# Row 116
$ws.Range("H116").Value = 1853.9474
$ws.Range("I116").Value = 1706.7333
$ws.Range("K116").Value = 1706.7333
$ws.Range("M116").Value = 587.2666999999999

# Row 122
$ws.Range("H122").Value = 4479.4
$ws.Range("I122").Value = 3849.6667
$ws.Range("J122").Value = 6998.3335
$ws.Range("K122").Value = 11549.0001
$ws.Range("L122").Value = 20995.0005
$ws.Range("M122").Value = -9099.000100000001
$ws.Range("N122").Value = -25895.0005

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1853.9474
$ws.Range("I3").Value = 1706.7333
$ws.Range("K3").Value = 1706.7333
$ws.Range("M3").Value = -1592.7333

# Row 4
$ws.Range("H4").Value = 183.75
$ws.Range("J4").Value = 315
$ws.Range("L4").Value = 315
$ws.Range("N4").Value = -545

# Row 11
$ws.Range("H11").Value = 1396.5
$ws.Range("I11").Value = 443
$ws.Range("J11").Value = 2350
$ws.Range("K11").Value = 443
$ws.Range("L11").Value = 2350
$ws.Range("M11").Value = -303
$ws.Range("N11").Value = -2630

# Row 105
$ws.Range("H105").Value = 2647.6667
$ws.Range("I105").Value = 2579.2727
$ws.Range("K105").Value = 2579.2727
$ws.Range("M105").Value = -832.2727

# Row 134
$ws.Range("H134").Value = 3007.1304
$ws.Range("I134").Value = 2373.7778
$ws.Range("K134").Value = 7121.3334
$ws.Range("M134").Value = -4586.3334

# Row 138
$ws.Range("H138").Value = 69499
$ws.Range("J138").Value = 69499
$ws.Range("L138").Value = 69499
$ws.Range("N138").Value = -79779

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 3399.6667
$ws.Range("I6").Value = 4524.5
$ws.Range("J6").Value = 1150
$ws.Range("K6").Value = 4524.5
$ws.Range("L6").Value = 1150
$ws.Range("M6").Value = -4411.5
$ws.Range("N6").Value = -1376

# Row 31
$ws.Range("H31").Value = 272132.1
$ws.Range("I31").Value = 371688.47
$ws.Range("J31").Value = 3329.8
$ws.Range("K31").Value = 371688.47
$ws.Range("L31").Value = 3329.8
$ws.Range("M31").Value = -371393.47
$ws.Range("N31").Value = -3919.8

# Row 34
$ws.Range("H34").Value = 272132.1
$ws.Range("I34").Value = 371688.47
$ws.Range("J34").Value = 3329.8
$ws.Range("K34").Value = 371688.47
$ws.Range("L34").Value = 3329.8
$ws.Range("M34").Value = -371486.47
$ws.Range("N34").Value = -3733.8

# Row 64
$ws.Range("H64").Value = 42088
$ws.Range("J64").Value = 42088
$ws.Range("L64").Value = 42088
$ws.Range("N64").Value = -42584

# Row 67
$ws.Range("H67").Value = 42088
$ws.Range("J67").Value = 42088
$ws.Range("L67").Value = 42088
$ws.Range("N67").Value = -43804

# Row 68
$ws.Range("H68").Value = 36826.43
$ws.Range("J68").Value = 36826.43
$ws.Range("L68").Value = 36826.43
$ws.Range("N68").Value = -38324.43

# Row 71
$ws.Range("H71").Value = 36826.43
$ws.Range("J71").Value = 36826.43
$ws.Range("L71").Value = 110479.29
$ws.Range("N71").Value = -117967.29

# Row 122
$ws.Range("H122").Value = 4530.25
$ws.Range("I122").Value = 4330
$ws.Range("K122").Value = 12990
$ws.Range("M122").Value = -10540

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("J131").Value = 1623.2858
$ws.Range("L131").Value = 4869.857400000001
$ws.Range("N131").Value = -14949.8574

# Row 140
$ws.Range("H140").Value = 29425188
$ws.Range("I140").Value = 38477370
$ws.Range("J140").Value = 5600
$ws.Range("K140").Value = 115432110
$ws.Range("L140").Value = 16800
$ws.Range("M140").Value = -115426930
$ws.Range("N140").Value = -27160

$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 16183.8
$ws.Range("I12").Value = 7973.3335
$ws.Range("J12").Value = 28499.5
$ws.Range("K12").Value = 7973.3335
$ws.Range("L12").Value = 28499.5
$ws.Range("M12").Value = -7833.3335
$ws.Range("N12").Value = -28779.5

# Row 113
$ws.Range("H113").Value = 6613.909
$ws.Range("I113").Value = 10491
$ws.Range("K113").Value = 10491
$ws.Range("M113").Value = -8321

# Row 132
$ws.Range("H132").Value = 33932.85
$ws.Range("I132").Value = 36771.71
$ws.Range("K132").Value = 110315.13
$ws.Range("M132").Value = -107785.13

$ws = $wb.Worksheets.Item("LTW")
# Row 38
$ws.Range("H38").Value = 58248.25
$ws.Range("J38").Value = 68999.5
$ws.Range("L38").Value = 68999.5
$ws.Range("N38").Value = -69819.5

# Row 68
$ws.Range("H68").Value = 8860.315000000001
$ws.Range("I68").Value = 9670.412
$ws.Range("J68").Value = 1974.5
$ws.Range("K68").Value = 9670.412
$ws.Range("L68").Value = 1974.5
$ws.Range("M68").Value = -8921.412
$ws.Range("N68").Value = -3472.5

# Row 71
$ws.Range("H71").Value = 8860.315000000001
$ws.Range("I71").Value = 9670.412
$ws.Range("J71").Value = 1974.5
$ws.Range("K71").Value = 48352.06
$ws.Range("L71").Value = 9872.5
$ws.Range("M71").Value = -44608.06
$ws.Range("N71").Value = -17360.5

# Row 93
$ws.Range("H93").Value = 111114820
$ws.Range("I93").Value = 3579.8
$ws.Range("J93").Value = 250003860
$ws.Range("K93").Value = 3579.8
$ws.Range("L93").Value = 250003860
$ws.Range("M93").Value = -2331.8
$ws.Range("N93").Value = -250006356

# Row 122
$ws.Range("H122").Value = 3192.2354
$ws.Range("I122").Value = 3178.9092
$ws.Range("K122").Value = 9536.7276
$ws.Range("M122").Value = -7086.7276

# Row 136
$ws.Range("H136").Value = 1274.56
$ws.Range("I136").Value = 1145.5264
$ws.Range("K136").Value = 3436.5792
$ws.Range("M136").Value = -886.5792000000001
